$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 67, shifting existing rows 67-98 down to 68-99.
$ws.Rows.Item(67).Insert()

# Populate the newly inserted row 67 with the new weekly record.
$ws.Range("A67").Value = 1
$ws.Range("B67").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C67").Value = "Arica y Parinacota"
$ws.Range("D67").Value = 45146
$ws.Range("E67").Value = 15
$ws.Range("F67").Value = 100112009
$ws.Range("G67").Value = "Acelga"
$ws.Range("H67").Value = "Sin especificar"
$ws.Range("I67").Value = "Primera"
$ws.Range("J67").Value = 350
$ws.Range("K67").Value = 1000
$ws.Range("L67").Value = 1200
$ws.Range("M67").Value = 1114
$ws.Range("N67").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O67").Value = "Región de Arica y Parinacota"
$ws.Range("P67").Value = 371
$ws.Range("Q67").Value = 3
$ws.Range("R67").Value = "Hortaliza"
